$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest coinranking.com snapshot: (cell -> new text) pairs plus a flag for
# cells whose text happens to look like a plain number (e.g. '43.95'). Those
# are written with a temporary Text number format so Excel's COM layer keeps
# them as literal strings instead of silently coercing them to doubles; the
# cell style is then restored to Normal so no stray formatting is left behind.
$updates = @(
    @{ Cell = 'D2'; Value = '36.505.13'; NumericLooking = $false }
    @{ Cell = 'E2'; Value = '  +3.40%  '; NumericLooking = $false }
    @{ Cell = 'D3'; Value = '1.917.84'; NumericLooking = $false }
    @{ Cell = 'E3'; Value = '  +1.72%  '; NumericLooking = $false }
    @{ Cell = 'E4'; Value = '  -0.05%  '; NumericLooking = $false }
    @{ Cell = 'D5'; Value = '249.13'; NumericLooking = $true }
    @{ Cell = 'E5'; Value = '  +1.40%  '; NumericLooking = $false }
    @{ Cell = 'E6'; Value = '  +0.76%  '; NumericLooking = $false }
    @{ Cell = 'D7'; Value = '1.00'; NumericLooking = $true }
    @{ Cell = 'E7'; Value = '  -0.03%  '; NumericLooking = $false }
    @{ Cell = 'D8'; Value = '43.95'; NumericLooking = $true }
    @{ Cell = 'E8'; Value = '  +1.10%  '; NumericLooking = $false }
    @{ Cell = 'D9'; Value = '58.30'; NumericLooking = $true }
    @{ Cell = 'E9'; Value = '  +8.72%  '; NumericLooking = $false }
    @{ Cell = 'E10'; Value = '  +3.15%  '; NumericLooking = $false }
    @{ Cell = 'E11'; Value = '  +3.13%  '; NumericLooking = $false }
    @{ Cell = 'D12'; Value = '0.0993'; NumericLooking = $true }
    @{ Cell = 'E12'; Value = '  +2.26%  '; NumericLooking = $false }
    @{ Cell = 'D13'; Value = '14.45'; NumericLooking = $true }
    @{ Cell = 'D14'; Value = '0.802'; NumericLooking = $true }
    @{ Cell = 'E14'; Value = '  +6.28%  '; NumericLooking = $false }
    @{ Cell = 'D15'; Value = '2.193.83'; NumericLooking = $false }
    @{ Cell = 'E15'; Value = '  +1.57%  '; NumericLooking = $false }
    @{ Cell = 'E16'; Value = '  +4.71%  '; NumericLooking = $false }
    @{ Cell = 'D17'; Value = '1.909.04'; NumericLooking = $false }
    @{ Cell = 'E17'; Value = '  +0.67%  '; NumericLooking = $false }
    @{ Cell = 'D18'; Value = '36.425.70'; NumericLooking = $false }
    @{ Cell = 'E18'; Value = '  +2.77%  '; NumericLooking = $false }
    @{ Cell = 'D19'; Value = '74.36'; NumericLooking = $true }
    @{ Cell = 'E19'; Value = '  +1.98%  '; NumericLooking = $false }
    @{ Cell = 'E20'; Value = '  +3.55%  '; NumericLooking = $false }
    @{ Cell = 'D21'; Value = '252.18'; NumericLooking = $true }
    @{ Cell = 'E21'; Value = '  +3.25%  '; NumericLooking = $false }
    @{ Cell = 'D22'; Value = '13.19'; NumericLooking = $true }
    @{ Cell = 'E22'; Value = '  +3.31%  '; NumericLooking = $false }
    @{ Cell = 'D23'; Value = '5.18'; NumericLooking = $true }
    @{ Cell = 'E23'; Value = '  +4.78%  '; NumericLooking = $false }
    @{ Cell = 'E24'; Value = '  +1.95%  '; NumericLooking = $false }
    @{ Cell = 'E25'; Value = '  +0.01%  '; NumericLooking = $false }
    @{ Cell = 'D26'; Value = '2.20'; NumericLooking = $true }
    @{ Cell = 'E26'; Value = '  +2.76%  '; NumericLooking = $false }
    @{ Cell = 'E27'; Value = '  +1.13%  '; NumericLooking = $false }
    @{ Cell = 'D28'; Value = '8.76'; NumericLooking = $true }
    @{ Cell = 'E28'; Value = '  +3.18%  '; NumericLooking = $false }
    @{ Cell = 'D29'; Value = '18.83'; NumericLooking = $true }
    @{ Cell = 'E29'; Value = '  +2.98%  '; NumericLooking = $false }
    @{ Cell = 'E30'; Value = '  +1.45%  '; NumericLooking = $false }
    @{ Cell = 'E31'; Value = '  +6.11%  '; NumericLooking = $false }
    @{ Cell = 'E32'; Value = '  +4.07%  '; NumericLooking = $false }
    @{ Cell = 'D33'; Value = '1.95'; NumericLooking = $true }
    @{ Cell = 'E33'; Value = '  +6.49%  '; NumericLooking = $false }
    @{ Cell = 'D34'; Value = '4.33'; NumericLooking = $true }
    @{ Cell = 'E34'; Value = '  +4.44%  '; NumericLooking = $false }
    @{ Cell = 'E35'; Value = '  +0.00%  '; NumericLooking = $false }
    @{ Cell = 'D36'; Value = '0.0849'; NumericLooking = $true }
    @{ Cell = 'E36'; Value = '  +22.39%  '; NumericLooking = $false }
    @{ Cell = 'E37'; Value = '  -15.03%  '; NumericLooking = $false }
    @{ Cell = 'D38'; Value = '0.859'; NumericLooking = $true }
    @{ Cell = 'E38'; Value = '  +1.23%  '; NumericLooking = $false }
    @{ Cell = 'D39'; Value = '2.01'; NumericLooking = $true }
    @{ Cell = 'E39'; Value = '  +2.44%  '; NumericLooking = $false }
    @{ Cell = 'D40'; Value = '106.52'; NumericLooking = $true }
    @{ Cell = 'E40'; Value = '  +10.37%  '; NumericLooking = $false }
    @{ Cell = 'B41'; Value = 'Gas'; NumericLooking = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'; NumericLooking = $false }
    @{ Cell = 'D41'; Value = '16.35'; NumericLooking = $true }
    @{ Cell = 'E41'; Value = '  +32.78%  '; NumericLooking = $false }
    @{ Cell = 'B42'; Value = 'VeChain'; NumericLooking = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; NumericLooking = $false }
    @{ Cell = 'D42'; Value = '0.0229'; NumericLooking = $true }
    @{ Cell = 'E42'; Value = '  +3.98%  '; NumericLooking = $false }
    @{ Cell = 'D43'; Value = '17.02'; NumericLooking = $true }
    @{ Cell = 'E43'; Value = '  -1.16%  '; NumericLooking = $false }
    @{ Cell = 'E44'; Value = '  +2.85%  '; NumericLooking = $false }
    @{ Cell = 'D45'; Value = '1.343.35'; NumericLooking = $false }
    @{ Cell = 'E45'; Value = '  +3.39%  '; NumericLooking = $false }
    @{ Cell = 'D46'; Value = '2.36'; NumericLooking = $true }
    @{ Cell = 'E46'; Value = '  +1.85%  '; NumericLooking = $false }
    @{ Cell = 'D47'; Value = '0.0806'; NumericLooking = $true }
    @{ Cell = 'E47'; Value = '  +1.11%  '; NumericLooking = $false }
    @{ Cell = 'E48'; Value = '  +2.60%  '; NumericLooking = $false }
    @{ Cell = 'D49'; Value = '2.79'; NumericLooking = $true }
    @{ Cell = 'E49'; Value = '  +2.14%  '; NumericLooking = $false }
    @{ Cell = 'E50'; Value = '  +3.06%  '; NumericLooking = $false }
    @{ Cell = 'D51'; Value = '2.092.53'; NumericLooking = $false }
    @{ Cell = 'E51'; Value = '  +1.15%  '; NumericLooking = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.NumericLooking) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
